# Daily "days remaining" countdown update.
# For every data row, decrement the "剩余" (remaining) value in column E
# by 1 -- except rows whose "开始时间" (start date, column F) is not a
# well-formed 8-digit yyyyMMdd value (malformed/typo'd dates are skipped,
# i.e. left completely untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
if ($lastRow -lt 2) {
    $lastRow = 99
}

for ($r = 2; $r -le $lastRow; $r++) {
    $dateCell = $ws.Cells.Item($r, 6)
    $dateVal = $dateCell.Value2
    if ($dateVal -eq $null -or $dateVal -eq "") {
        continue
    }
    $dateStr = [string]$dateVal

    # Skip rows whose start-date (column F) is not a well-formed 8-digit
    # yyyyMMdd value (e.g. malformed/typo'd dates) - these are left untouched.
    if ($dateStr.Length -ne 8) {
        continue
    }

    $remCell = $ws.Cells.Item($r, 5)
    $remVal = $remCell.Value2
    if ($remVal -eq $null -or $remVal -eq "") {
        continue
    }

    $remCell.Value = [double]$remVal - 1
}
